# Apply updates to the Jogos_do_Dia Betfair Back/Lay workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric odds updates ---
$ws.Range("K2").Value = 6.8

$ws.Range("F4").Value = 3.05
$ws.Range("G4").Value = 11.5
$ws.Range("P4").Value = 2.3

$ws.Range("F5").Value = 1.55
$ws.Range("H5").Value = 2.22
$ws.Range("I5").Value = 2.9

$ws.Range("Q6").Value = 1.44

$ws.Range("I7").Value = 3.15

$ws.Range("G10").Value = 1.81
$ws.Range("H10").Value = 2.28
$ws.Range("J10").Value = 2.22

# --- SnapshotTS (BH column) updates for rows 2-11, keep as text ---
$newTimestamp = "2026-02-25 05:24:15"
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("BH$r").NumberFormat = "@"
    $ws.Range("BH$r").Value = $newTimestamp
}
